$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bachelor of Engineering paragraph: append " - Score: 82%" after
#    "Completed in 2023" as four separate runs (space / "- Score: " / "82" / "%").
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Completed in 2023", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) # wdCollapseEnd
$rng.InsertAfter(" ") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter([char]0x2013 + " Score: ") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("82") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("%") | Out-Null

# ---------------------------------------------------------------------------
# 2) Higher Secondary score: "68.50%" -> "68%"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("68.50%", $true, $false, $false, $false, $false, $true, 1, $false, "68%", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Date of Birth: "November 14, 2002" -> "July 15, 2002" split across runs
#    (" " / "July 15" / ", 2002")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Date of Birth:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraRange = $rng.Paragraphs(1).Range
$paraRange.Collapse(1) # wdCollapseStart

$dobXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Date of Birth:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>July 15</w:t></w:r><w:r><w:t>, 2002</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$paraRange.InsertXML($dobXml) | Out-Null

# ---------------------------------------------------------------------------
# 4) Location: " Manjarapatti, Pennagaram, Dharmapuri, Tamil Nadu" ->
#    " Manjarapatti, Pennagaram(Tk), Dharmapuri(Dt), Tamil Nadu-636810"
#    with proofing-error markers around the place names, as Word's spell
#    checker would insert them.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Location:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraRange = $rng.Paragraphs(1).Range
$paraRange.Collapse(1) # wdCollapseStart

$locXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Location:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Manjarapatti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Pennagaram</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Tk)</w:t></w:r><w:r><w:t>, Dharmapuri</w:t></w:r><w:r><w:t>(Dt)</w:t></w:r><w:r><w:t>, Tamil Nadu</w:t></w:r><w:r><w:t>-636810</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$paraRange.InsertXML($locXml) | Out-Null
